$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# hunk 0: ALC!row43
$ws.Cells.Item(43, 8).Value = 1748.1904
$ws.Cells.Item(43, 9).Value = 3105
$ws.Cells.Item(43, 10).Value = 1205.4667
$ws.Cells.Item(43, 11).Value = 3105
$ws.Cells.Item(43, 12).Value = 1205.4667
$ws.Cells.Item(43, 13).Value = -3036
$ws.Cells.Item(43, 14).Value = -1343.4667

# hunk 1: ALC!row62
$ws.Cells.Item(62, 8).Value = 348949.22
$ws.Cells.Item(62, 9).Value = 529916.75
$ws.Cells.Item(62, 10).Value = 5110.9
$ws.Cells.Item(62, 11).Value = 529916.75
$ws.Cells.Item(62, 12).Value = 5110.9
$ws.Cells.Item(62, 13).Value = -529292.75
$ws.Cells.Item(62, 14).Value = -6358.9

# hunk 2: ALC!row64
$ws.Cells.Item(64, 8).Value = 3491.96
$ws.Cells.Item(64, 9).Value = 3408.2917
$ws.Cells.Item(64, 11).Value = 3408.2917
$ws.Cells.Item(64, 13).Value = -3160.2917

# hunk 3: ALC!row65
$ws.Cells.Item(65, 8).Value = 348949.22
$ws.Cells.Item(65, 9).Value = 529916.75
$ws.Cells.Item(65, 10).Value = 5110.9
$ws.Cells.Item(65, 11).Value = 2649583.75
$ws.Cells.Item(65, 12).Value = 25554.5
$ws.Cells.Item(65, 13).Value = -2646463.75
$ws.Cells.Item(65, 14).Value = -31794.5

# hunk 4: ALC!row67
$ws.Cells.Item(67, 8).Value = 3491.96
$ws.Cells.Item(67, 9).Value = 3408.2917
$ws.Cells.Item(67, 11).Value = 3408.2917
$ws.Cells.Item(67, 13).Value = -2550.2917

# hunk 5: ALC!row105
$ws.Cells.Item(105, 8).Value = 39800
$ws.Cells.Item(105, 10).Value = 39800
$ws.Cells.Item(105, 12).Value = 39800
$ws.Cells.Item(105, 14).Value = -46788

# hunk 6: ALC!row113
$ws.Cells.Item(113, 8).Value = 3869.3125
$ws.Cells.Item(113, 9).Value = 3475.25
$ws.Cells.Item(113, 11).Value = 3475.25
$ws.Cells.Item(113, 13).Value = -221.25

# hunk 7: ALC!row125
$ws.Cells.Item(125, 8).Value = 936.9375
$ws.Cells.Item(125, 9).Value = 926.5
$ws.Cells.Item(125, 10).Value = 954.3333
$ws.Cells.Item(125, 11).Value = 8338.5
$ws.Cells.Item(125, 12).Value = 8588.9997
$ws.Cells.Item(125, 13).Value = -5878.5
$ws.Cells.Item(125, 14).Value = -13508.9997

# hunk 8: ALC!row129
$ws.Cells.Item(129, 8).Value = 1094.2759
$ws.Cells.Item(129, 9).Value = 404.85715
$ws.Cells.Item(129, 10).Value = 1313.6364
$ws.Cells.Item(129, 11).Value = 1214.57145
$ws.Cells.Item(129, 12).Value = 3940.9092
$ws.Cells.Item(129, 13).Value = 3785.42855
$ws.Cells.Item(129, 14).Value = -13940.9092

# hunk 9: ALC!row138
$ws.Cells.Item(138, 8).Value = 3135.2842
$ws.Cells.Item(138, 9).Value = 1026.3513
$ws.Cells.Item(138, 10).Value = 4480.6377
$ws.Cells.Item(138, 11).Value = 3079.0539
$ws.Cells.Item(138, 12).Value = 13441.9131
$ws.Cells.Item(138, 13).Value = 2060.9461
$ws.Cells.Item(138, 14).Value = -23721.9131

$ws = $wb.Worksheets.Item("ARM")
# hunk 10: ARM!row2
$ws.Cells.Item(2, 8).Value = 2661.4736
$ws.Cells.Item(2, 9).Value = 2273.4614
$ws.Cells.Item(2, 10).Value = 3502.1667
$ws.Cells.Item(2, 11).Value = 2273.4614
$ws.Cells.Item(2, 12).Value = 3502.1667
$ws.Cells.Item(2, 13).Value = -2160.4614
$ws.Cells.Item(2, 14).Value = -3728.1667

# hunk 11: ARM!row32
$ws.Cells.Item(32, 8).Value = 2521.06
$ws.Cells.Item(32, 9).Value = 2521.06
$ws.Cells.Item(32, 11).Value = 2521.06
$ws.Cells.Item(32, 13).Value = -2234.06

# hunk 12: ARM!row45
$ws.Cells.Item(45, 8).Value = 972.0345
$ws.Cells.Item(45, 9).Value = 898.85
$ws.Cells.Item(45, 10).Value = 1134.6666
$ws.Cells.Item(45, 11).Value = 898.85
$ws.Cells.Item(45, 12).Value = 1134.6666
$ws.Cells.Item(45, 13).Value = -521.85
$ws.Cells.Item(45, 14).Value = -1888.6666

# hunk 13: ARM!row116
$ws.Cells.Item(116, 8).Value = 2661.4736
$ws.Cells.Item(116, 9).Value = 2273.4614
$ws.Cells.Item(116, 10).Value = 3502.1667
$ws.Cells.Item(116, 11).Value = 2273.4614
$ws.Cells.Item(116, 12).Value = 3502.1667
$ws.Cells.Item(116, 13).Value = 20.53859999999986
$ws.Cells.Item(116, 14).Value = -8090.1667

# hunk 14: ARM!row122
$ws.Cells.Item(122, 8).Value = 1350
$ws.Cells.Item(122, 9).Value = 1200
$ws.Cells.Item(122, 11).Value = 3600
$ws.Cells.Item(122, 13).Value = -1150

$ws = $wb.Worksheets.Item("BSM")
# hunk 15: BSM!row3
$ws.Cells.Item(3, 8).Value = 2661.4736
$ws.Cells.Item(3, 9).Value = 2273.4614
$ws.Cells.Item(3, 10).Value = 3502.1667
$ws.Cells.Item(3, 11).Value = 2273.4614
$ws.Cells.Item(3, 12).Value = 3502.1667
$ws.Cells.Item(3, 13).Value = -2159.4614
$ws.Cells.Item(3, 14).Value = -3730.1667

# hunk 16: BSM!row20
$ws.Cells.Item(20, 8).Value = 33359792
$ws.Cells.Item(20, 9).Value = 35272.895
$ws.Cells.Item(20, 10).Value = 90920330
$ws.Cells.Item(20, 11).Value = 35272.895
$ws.Cells.Item(20, 12).Value = 90920330
$ws.Cells.Item(20, 13).Value = -35025.895
$ws.Cells.Item(20, 14).Value = -90920824

$ws = $wb.Worksheets.Item("CRP")
# hunk 17: CRP!row122
$ws.Cells.Item(122, 8).Value = 1599.8
$ws.Cells.Item(122, 9).Value = 1733
$ws.Cells.Item(122, 10).Value = 1400
$ws.Cells.Item(122, 11).Value = 5199
$ws.Cells.Item(122, 12).Value = 4200
$ws.Cells.Item(122, 13).Value = -2749
$ws.Cells.Item(122, 14).Value = -9100

$ws = $wb.Worksheets.Item("CUL")
# hunk 18: CUL!row113
$ws.Cells.Item(113, 8).Value = 575.7959
$ws.Cells.Item(113, 10).Value = 550.04877
$ws.Cells.Item(113, 12).Value = 1650.14631
$ws.Cells.Item(113, 14).Value = -5990.14631

# hunk 19: CUL!row140
$ws.Cells.Item(140, 8).Value = 127666.25
$ws.Cells.Item(140, 9).Value = 202872.67
$ws.Cells.Item(140, 10).Value = 2322.2222
$ws.Cells.Item(140, 11).Value = 608618.01
$ws.Cells.Item(140, 12).Value = 6966.6666
$ws.Cells.Item(140, 13).Value = -603438.01
$ws.Cells.Item(140, 14).Value = -17326.6666

$ws = $wb.Worksheets.Item("GSM")
# hunk 20: GSM!row70
$ws.Cells.Item(70, 8).Value = 4639.1665
$ws.Cells.Item(70, 9).Value = 4112.4287
$ws.Cells.Item(70, 10).Value = 4974.364
$ws.Cells.Item(70, 11).Value = 4112.4287
$ws.Cells.Item(70, 12).Value = 4974.364
$ws.Cells.Item(70, 13).Value = -3842.4287
$ws.Cells.Item(70, 14).Value = -5514.364

# hunk 21: GSM!row73
$ws.Cells.Item(73, 8).Value = 4639.1665
$ws.Cells.Item(73, 9).Value = 4112.4287
$ws.Cells.Item(73, 10).Value = 4974.364
$ws.Cells.Item(73, 11).Value = 4112.4287
$ws.Cells.Item(73, 12).Value = 4974.364
$ws.Cells.Item(73, 13).Value = -3176.4287
$ws.Cells.Item(73, 14).Value = -6846.364

# hunk 22: GSM!row80
$ws.Cells.Item(80, 8).Value = 3971.7856
$ws.Cells.Item(80, 9).Value = 4117.0835
$ws.Cells.Item(80, 10).Value = 3100
$ws.Cells.Item(80, 11).Value = 4117.0835
$ws.Cells.Item(80, 12).Value = 3100
$ws.Cells.Item(80, 13).Value = -3119.0835
$ws.Cells.Item(80, 14).Value = -5096

# hunk 23: GSM!row83
$ws.Cells.Item(83, 8).Value = 3971.7856
$ws.Cells.Item(83, 9).Value = 4117.0835
$ws.Cells.Item(83, 10).Value = 3100
$ws.Cells.Item(83, 11).Value = 20585.4175
$ws.Cells.Item(83, 12).Value = 15500
$ws.Cells.Item(83, 13).Value = -15593.4175
$ws.Cells.Item(83, 14).Value = -25484

# hunk 24: GSM!row102
$ws.Cells.Item(102, 8).Value = 1314.05
$ws.Cells.Item(102, 9).Value = 1286.0667
$ws.Cells.Item(102, 10).Value = 1398
$ws.Cells.Item(102, 11).Value = 1286.0667
$ws.Cells.Item(102, 12).Value = 1398
$ws.Cells.Item(102, 13).Value = 335.9332999999999
$ws.Cells.Item(102, 14).Value = -4642

# hunk 25: GSM!row122
$ws.Cells.Item(122, 8).Value = 1138.2354
$ws.Cells.Item(122, 9).Value = 1175.8572
$ws.Cells.Item(122, 11).Value = 3527.5716
$ws.Cells.Item(122, 13).Value = -1077.5716

# hunk 26: GSM!row126
$ws.Cells.Item(126, 8).Value = 4362.4
$ws.Cells.Item(126, 9).Value = 4703
$ws.Cells.Item(126, 10).Value = 3000
$ws.Cells.Item(126, 11).Value = 14109
$ws.Cells.Item(126, 12).Value = 9000
$ws.Cells.Item(126, 13).Value = -11639
$ws.Cells.Item(126, 14).Value = -13940

# hunk 27: GSM!row132
$ws.Cells.Item(132, 8).Value = 2262.8
$ws.Cells.Item(132, 9).Value = 2088.5
$ws.Cells.Item(132, 10).Value = 2960
$ws.Cells.Item(132, 11).Value = 6265.5
$ws.Cells.Item(132, 12).Value = 8880
$ws.Cells.Item(132, 13).Value = -3735.5
$ws.Cells.Item(132, 14).Value = -13940

$ws = $wb.Worksheets.Item("LTW")
# hunk 28: LTW!row46
$ws.Cells.Item(46, 8).Value = 1162
$ws.Cells.Item(46, 9).Value = 1624.5
$ws.Cells.Item(46, 11).Value = 1624.5
$ws.Cells.Item(46, 13).Value = -1436.5

# hunk 29: LTW!row122
$ws.Cells.Item(122, 8).Value = 1959.6
$ws.Cells.Item(122, 9).Value = 1959.6
$ws.Cells.Item(122, 11).Value = 5878.799999999999
$ws.Cells.Item(122, 13).Value = -3428.799999999999

$ws = $wb.Worksheets.Item("WVR")
# hunk 30: WVR!row122
$ws.Cells.Item(122, 8).Value = 43479576
$ws.Cells.Item(122, 9).Value = 47620348
$ws.Cells.Item(122, 11).Value = 142861044
$ws.Cells.Item(122, 13).Value = -142858594

Write-Output "done"